$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.222.30'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.904.91'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3809'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07307'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9057'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08057'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("D15").Value = '1.795.08'
$ws.Range("E15").Value = '  -5.76%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008682'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '27.257.44'
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.128'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  +1.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.477'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.355'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.846'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.898'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09252'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8047'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05082'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.228'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.987'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.389'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.695'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.23%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.087'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.003'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.607'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4924'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  +1.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05961'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.37%  '
